$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point conversions (914400 EMU/in, 12700 EMU/pt):
#   off  x=498909   y=365760  -> 39.284173228346454 pt, 28.8 pt
#   ext cx=11194181 cy=923330 -> 881.4315748031496 pt, 72.7031496062992 pt
$shp = $s.Shapes.AddTextbox(1, 39.284173228346454, 28.8, 881.4315748031496, 72.7031496062992)
$shp.Name = "TextBox 1"

$shp.TextFrame.WordWrap = $true
$shp.Fill.Visible = $false

$tr = $shp.TextFrame.TextRange
$tr.Text = "이곳에 텍스트 입력"
$tr.LanguageID = "ko-KR"
$tr.ParagraphFormat.Alignment = 2
$tr.Font.Size = 54
$tr.Font.Bold = $true

# Resize-shape-to-fit-text (<a:spAutoFit/>) recomputed last so the
# rendered extent reflects the final font/bold/size settings above.
$shp.TextFrame.AutoSize = 1
